# Apply updated cryptos list values (price + 1h volume change) per row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '28.374.24'
$cell.Style = $origStyle
$ws.Range("E2").Value = '  +3.42%  '

$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.863.79'
$cell.Style = $origStyle
$ws.Range("E3").Value = '  +2.04%  '

$ws.Range("E4").Value = '  -0.58%  '

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '337.02'
$cell.Style = $origStyle
$ws.Range("E5").Value = '  +1.80%  '

$ws.Range("E6").Value = '  -0.56%  '

$cell = $ws.Range("D7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.4715'
$cell.Style = $origStyle
$ws.Range("E7").Value = '  +3.19%  '

$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.3967'
$cell.Style = $origStyle
$ws.Range("E8").Value = '  +3.87%  '

$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '47.55'
$cell.Style = $origStyle
$ws.Range("E9").Value = '  +2.24%  '

$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.08002'
$cell.Style = $origStyle
$ws.Range("E10").Value = '  +1.23%  '

$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.9943'
$cell.Style = $origStyle
$ws.Range("E11").Value = '  +2.84%  '

$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '21.91'
$cell.Style = $origStyle
$ws.Range("E12").Value = '  +4.15%  '

$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.021'
$cell.Style = $origStyle
$ws.Range("E13").Value = '  +2.66%  '

$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.867.36'
$cell.Style = $origStyle
$ws.Range("E14").Value = '  +0.94%  '

$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '7.243'
$cell.Style = $origStyle
$ws.Range("E15").Value = '  +2.92%  '

$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '90.32'
$cell.Style = $origStyle
$ws.Range("E16").Value = '  +2.55%  '

$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.Style = $origStyle
$ws.Range("E17").Value = '  -0.58%  '

$ws.Range("E18").Value = '  +0.92%  '

$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.06612'
$cell.Style = $origStyle
$ws.Range("E19").Value = '  -0.43%  '

$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '17.47'
$cell.Style = $origStyle
$ws.Range("E20").Value = '  +1.85%  '

$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.Style = $origStyle
$ws.Range("E21").Value = '  -0.45%  '

$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '28.390.58'
$cell.Style = $origStyle
$ws.Range("E22").Value = '  +3.51%  '

$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.458'
$cell.Style = $origStyle
$ws.Range("E23").Value = '  +2.45%  '

$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '11.01'
$cell.Style = $origStyle
$ws.Range("E24").Value = '  +2.19%  '

$ws.Range("E25").Value = '  -1.45%  '

$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.086.27'
$cell.Style = $origStyle
$ws.Range("E26").Value = '  +1.06%  '

$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '161.37'
$cell.Style = $origStyle
$ws.Range("E27").Value = '  +2.72%  '

$ws.Range("E28").Value = '  +1.74%  '

$ws.Range("E29").Value = '  +2.57%  '

$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.450'
$cell.Style = $origStyle
$ws.Range("E30").Value = '  +4.35%  '

$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '119.20'
$cell.Style = $origStyle
$ws.Range("E31").Value = '  +0.91%  '

$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.09525'
$cell.Style = $origStyle
$ws.Range("E32").Value = '  +2.60%  '

$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.9572'
$cell.Style = $origStyle
$ws.Range("E33").Value = '  +1.40%  '

$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.592'
$cell.Style = $origStyle
$ws.Range("E34").Value = '  +0.54%  '

$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.373'
$cell.Style = $origStyle
$ws.Range("E35").Value = '  +4.70%  '

$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.344'
$cell.Style = $origStyle
$ws.Range("E36").Value = '  +2.13%  '

$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.06131'
$cell.Style = $origStyle
$ws.Range("E37").Value = '  +3.48%  '

$ws.Range("E38").Value = '  +2.65%  '

$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '8.267'
$cell.Style = $origStyle
$ws.Range("E39").Value = '  +3.61%  '

$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.177'
$cell.Style = $origStyle
$ws.Range("E40").Value = '  +1.74%  '

$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.5911'
$cell.Style = $origStyle
$ws.Range("E41").Value = '  +2.29%  '

$ws.Range("E42").Value = '  -0.52%  '

$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.1873'
$cell.Style = $origStyle
$ws.Range("E43").Value = '  +2.25%  '

$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '10.26'
$cell.Style = $origStyle
$ws.Range("E44").Value = '  +2.68%  '

$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.270'
$cell.Style = $origStyle
$ws.Range("E45").Value = '  -0.56%  '

$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.07537'
$cell.Style = $origStyle
$ws.Range("E46").Value = '  +13.49%  '

$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.5532'
$cell.Style = $origStyle
$ws.Range("E47").Value = '  +1.06%  '

$ws.Range("E48").Value = '  +0.77%  '

$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.940'
$cell.Style = $origStyle
$ws.Range("E49").Value = '  +4.11%  '

$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.056'
$cell.Style = $origStyle
$ws.Range("E50").Value = '  +12.84%  '

$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '111.88'
$cell.Style = $origStyle
$ws.Range("E51").Value = '  +1.91%  '
